# This edit refreshes the "Metrics" sheet's raw input numbers (B2:B13).
# The "today" sheet's B11:B22 read these via =Metrics!Bn formulas, and its
# E/F columns derive from those in turn, so once Metrics is updated and
# Excel recalculates, today's cached formula values follow automatically -
# no need to touch them (or the volatile =TODAY()-1 in today!A1) directly.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Activate()
$metrics.Range("B2").Value  = 145609.93
$metrics.Range("B3").Value  = 125037.34000000001
$metrics.Range("B4").Value  = 44822.09
$metrics.Range("B5").Value  = 5968
$metrics.Range("B6").Value  = 5348317.040000001
$metrics.Range("B7").Value  = 4525390.3000000007
$metrics.Range("B8").Value  = 1576778.9700000004
$metrics.Range("B9").Value  = 208675
$metrics.Range("B10").Value = 33813698.029999994
$metrics.Range("B11").Value = 31800665.460000001
$metrics.Range("B12").Value = 11858501.009999994
$metrics.Range("B13").Value = 1306305

# Move the Metrics sheet's own selection/active cell to G9, matching the
# recorded UI state after the edit.
$metrics.Range("G9").Select()

# The "today" sheet is the sheet that was active/tab-selected in the
# workbook, and stays so: activate it again and move its selection to E6.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E6").Select()
